# Apply trade #7 close to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet -----------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.14     # Total P&L %
$summary.Range("B6").Value = 7         # Total Trades
$summary.Range("B9").Value = 28.57     # Win Rate %

# --- Strategy Status sheet ---------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 7          # Trades (MarketMaking row)
$status.Range("G4").Value = 28.57      # Win Rate % (MarketMaking row)

# --- New trade row data (trade #7) -------------------------------------
$tradeNum    = 7
$tradeDate   = "2026-02-17"
$tradeTime   = "07:57:56"
$strategy    = "MarketMaking"
$side        = "DOWN"
$entryPrice  = 0.8
$exitPrice   = 0.8
$status7     = "CLOSED"
$pnlPct      = 0
$pnlDollar   = 0
$capAfter    = 99.95
$entrySlip   = 0
$exitSlip    = 0
$confidence  = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason  = "early_exit"
$duration    = 0.13

function Set-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value  = $tradeNum

    # Date/Time columns look numeric to Excel's auto-detection, so they'd
    # otherwise be stored as date serials. Force them to remain plain text
    # (matching columns B/C of the existing rows), then drop the leftover
    # "Text" number-format stamp so the cell style matches the other rows.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $tradeDate
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $tradeTime
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value  = $strategy
    $ws.Cells.Item($row, 5).Value  = $side
    $ws.Cells.Item($row, 6).Value  = $entryPrice
    $ws.Cells.Item($row, 7).Value  = $exitPrice
    $ws.Cells.Item($row, 8).Value  = $status7
    $ws.Cells.Item($row, 9).Value  = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capAfter
    $ws.Cells.Item($row, 12).Value = $entrySlip
    $ws.Cells.Item($row, 13).Value = $exitSlip
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $duration
}

# --- All Trades sheet ---------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Set-TradeRow $allTrades 8

# --- MarketMaking sheet --------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
Set-TradeRow $mm 8
